$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Cells.Item(2, 7).Value = 0.5023523333333334
$ws.Cells.Item(2, 8).Value = 1.507057
$ws.Cells.Item(2, 9).Value = 0.06515888850144765
$ws.Cells.Item(2, 10).Value = 0.06515888850144765
$ws.Cells.Item(2, 13).Value = 1.826566
$ws.Cells.Item(2, 14).Value = 5.479698
$ws.Cells.Item(2, 15).Value = 0.02795372904983374
$ws.Cells.Item(2, 16).Value = 0.02795372904983374
$ws.Cells.Item(2, 17).Value = 0.9175796920873334
$ws.Cells.Item(2, 18).Value = 8.258217228786
$ws.Cells.Item(2, 19).Value = 0.001821433914357795
$ws.Cells.Item(2, 20).Value = 0.001821433914357795

# Row 3 updates
$ws.Cells.Item(3, 7).Value = 0.5023523333333334
$ws.Cells.Item(3, 8).Value = 1.507057
$ws.Cells.Item(3, 9).Value = 0.06515888850144765
$ws.Cells.Item(3, 10).Value = 0.06515888850144765
$ws.Cells.Item(3, 13).Value = 44.29005966666667
$ws.Cells.Item(3, 15).Value = 0.6778141756295529
$ws.Cells.Item(3, 16).Value = 0.6778141756295529
$ws.Cells.Item(3, 17).Value = 22.24921481702256
$ws.Cells.Item(3, 18).Value = 200.242933353203
$ws.Cells.Item(3, 19).Value = 0.04416561829454669
$ws.Cells.Item(3, 20).Value = 0.04416561829454669

# Row 4 updates
$ws.Cells.Item(4, 7).Value = 0.5023523333333334
$ws.Cells.Item(4, 8).Value = 1.507057
$ws.Cells.Item(4, 9).Value = 0.06515888850144765
$ws.Cells.Item(4, 10).Value = 0.06515888850144765
$ws.Cells.Item(4, 13).Value = 19.10886933333333
$ws.Cells.Item(4, 14).Value = 57.326608
$ws.Cells.Item(4, 15).Value = 0.2924417490485847
$ws.Cells.Item(4, 16).Value = 0.2924417490485847
$ws.Cells.Item(4, 17).Value = 9.599385096961779
$ws.Cells.Item(4, 18).Value = 86.394465872656
$ws.Cells.Item(4, 19).Value = 0.01905517931942506
$ws.Cells.Item(4, 20).Value = 0.01905517931942506

# Row 5 updates
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 4).Value = "Resolving-Mac"
$ws.Cells.Item(5, 7).Value = 0.5023523333333334
$ws.Cells.Item(5, 8).Value = 1.507057
$ws.Cells.Item(5, 9).Value = 0.06515888850144765
$ws.Cells.Item(5, 10).Value = 0.06515888850144765
$ws.Cells.Item(5, 11).Value = 2.0
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.1169856666666667
$ws.Cells.Item(5, 14).Value = 0.350957
$ws.Cells.Item(5, 15).Value = 0.001790346272028586
$ws.Cells.Item(5, 16).Value = 0.001790346272028586
$ws.Cells.Item(5, 17).Value = 0.05876802261655556
$ws.Cells.Item(5, 18).Value = 0.5289122035489999
$ws.Cells.Item(5, 19).Value = 0.0001166569731180931
$ws.Cells.Item(5, 20).Value = 0.0001166569731180931

# Row 6 updates
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 9).Value = 0.07275905893716338
$ws.Cells.Item(6, 10).Value = 0.07275905893716339
$ws.Cells.Item(6, 13).Value = 1.826566
$ws.Cells.Item(6, 14).Value = 5.479698
$ws.Cells.Item(6, 15).Value = 0.02795372904983374
$ws.Cells.Item(6, 16).Value = 0.02795372904983374
$ws.Cells.Item(6, 17).Value = 1.024606718002
$ws.Cells.Item(6, 18).Value = 9.221460462018
$ws.Cells.Item(6, 19).Value = 0.002033887019450349
$ws.Cells.Item(6, 20).Value = 0.00203388701945035

# Row 7 updates
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 9).Value = 0.07275905893716338
$ws.Cells.Item(7, 10).Value = 0.07275905893716339
$ws.Cells.Item(7, 13).Value = 44.29005966666667
$ws.Cells.Item(7, 14).Value = 132.870179
$ws.Cells.Item(7, 15).Value = 0.6778141756295529
$ws.Cells.Item(7, 16).Value = 0.6778141756295529
$ws.Cells.Item(7, 17).Value = 24.84437609983767
$ws.Cells.Item(7, 18).Value = 223.599384898539
$ws.Cells.Item(7, 19).Value = 0.04931712155307545
$ws.Cells.Item(7, 20).Value = 0.04931712155307546

# Row 8 updates
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 4).Value = "MuSCs"
$ws.Cells.Item(8, 7).Value = 0.560947
$ws.Cells.Item(8, 8).Value = 1.682841
$ws.Cells.Item(8, 9).Value = 0.07275905893716338
$ws.Cells.Item(8, 10).Value = 0.07275905893716339
$ws.Cells.Item(8, 13).Value = 19.10886933333333
$ws.Cells.Item(8, 14).Value = 57.326608
$ws.Cells.Item(8, 15).Value = 0.2924417490485847
$ws.Cells.Item(8, 16).Value = 0.2924417490485847
$ws.Cells.Item(8, 17).Value = 10.71906292592533
$ws.Cells.Item(8, 18).Value = 96.471566333328
$ws.Cells.Item(8, 19).Value = 0.02127778645471312
$ws.Cells.Item(8, 20).Value = 0.02127778645471312

# Row 9 updates
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 4).Value = "Resolving-Mac"
$ws.Cells.Item(9, 7).Value = 0.560947
$ws.Cells.Item(9, 8).Value = 1.682841
$ws.Cells.Item(9, 9).Value = 0.07275905893716338
$ws.Cells.Item(9, 10).Value = 0.07275905893716339
$ws.Cells.Item(9, 11).Value = 2.0
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 0.1169856666666667
$ws.Cells.Item(9, 14).Value = 0.350957
$ws.Cells.Item(9, 15).Value = 0.001790346272028586
$ws.Cells.Item(9, 16).Value = 0.001790346272028586
$ws.Cells.Item(9, 17).Value = 0.06562275875966665
$ws.Cells.Item(9, 18).Value = 0.5906048288369999
$ws.Cells.Item(9, 19).Value = 0.0001302639099244587
$ws.Cells.Item(9, 20).Value = 0.0001302639099244587

# Row 10 updates
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 7).Value = 2.845667666666667
$ws.Cells.Item(10, 8).Value = 8.537003
$ws.Cells.Item(10, 9).Value = 0.3691045704399529
$ws.Cells.Item(10, 10).Value = 0.3691045704399529
$ws.Cells.Item(10, 13).Value = 1.826566
$ws.Cells.Item(10, 14).Value = 5.479698
$ws.Cells.Item(10, 15).Value = 0.02795372904983374
$ws.Cells.Item(10, 16).Value = 0.02795372904983374
$ws.Cells.Item(10, 17).Value = 5.197799807232667
$ws.Cells.Item(10, 18).Value = 46.780198265094
$ws.Cells.Item(10, 19).Value = 0.01031784915313371
$ws.Cells.Item(10, 20).Value = 0.01031784915313372

# Row 11 updates
$ws.Cells.Item(11, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 7).Value = 2.845667666666667
$ws.Cells.Item(11, 8).Value = 8.537003
$ws.Cells.Item(11, 9).Value = 0.3691045704399529
$ws.Cells.Item(11, 10).Value = 0.3691045704399529
$ws.Cells.Item(11, 13).Value = 44.29005966666667
$ws.Cells.Item(11, 14).Value = 132.870179
$ws.Cells.Item(11, 15).Value = 0.6778141756295529
$ws.Cells.Item(11, 16).Value = 0.6778141756295529
$ws.Cells.Item(11, 17).Value = 126.0347907481708
$ws.Cells.Item(11, 18).Value = 1134.313116733537
$ws.Cells.Item(11, 19).Value = 0.2501843101338569
$ws.Cells.Item(11, 20).Value = 0.2501843101338569

# Row 12 updates
$ws.Cells.Item(12, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(12, 4).Value = "MuSCs"
$ws.Cells.Item(12, 7).Value = 2.845667666666667
$ws.Cells.Item(12, 8).Value = 8.537003
$ws.Cells.Item(12, 9).Value = 0.3691045704399529
$ws.Cells.Item(12, 10).Value = 0.3691045704399529
$ws.Cells.Item(12, 13).Value = 19.10886933333333
$ws.Cells.Item(12, 14).Value = 57.326608
$ws.Cells.Item(12, 15).Value = 0.2924417490485847
$ws.Cells.Item(12, 16).Value = 0.2924417490485847
$ws.Cells.Item(12, 17).Value = 54.3774916084249
$ws.Cells.Item(12, 18).Value = 489.397424475824
$ws.Cells.Item(12, 19).Value = 0.1079415861612864
$ws.Cells.Item(12, 20).Value = 0.1079415861612864

# Row 13 updates
$ws.Cells.Item(13, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(13, 4).Value = "Resolving-Mac"
$ws.Cells.Item(13, 7).Value = 2.845667666666667
$ws.Cells.Item(13, 8).Value = 8.537003
$ws.Cells.Item(13, 9).Value = 0.3691045704399529
$ws.Cells.Item(13, 10).Value = 0.3691045704399529
$ws.Cells.Item(13, 11).Value = 2.0
$ws.Cells.Item(13, 12).Value = 0.6666666666666666
$ws.Cells.Item(13, 13).Value = 0.1169856666666667
$ws.Cells.Item(13, 14).Value = 0.350957
$ws.Cells.Item(13, 15).Value = 0.001790346272028586
$ws.Cells.Item(13, 16).Value = 0.001790346272028586
$ws.Cells.Item(13, 17).Value = 0.3329023290967777
$ws.Cells.Item(13, 18).Value = 2.996120961871
$ws.Cells.Item(13, 19).Value = 0.0006608249916758823
$ws.Cells.Item(13, 20).Value = 0.0006608249916758824

# Row 14 updates
$ws.Cells.Item(14, 1).Value = "MuSCs"
$ws.Cells.Item(14, 5).Value = 2.0
$ws.Cells.Item(14, 6).Value = 0.6666666666666666
$ws.Cells.Item(14, 7).Value = 0.11543
$ws.Cells.Item(14, 8).Value = 0.34629
$ws.Cells.Item(14, 9).Value = 0.01497214206175765
$ws.Cells.Item(14, 10).Value = 0.01497214206175765
$ws.Cells.Item(14, 13).Value = 1.826566
$ws.Cells.Item(14, 14).Value = 5.479698
$ws.Cells.Item(14, 15).Value = 0.02795372904983374
$ws.Cells.Item(14, 16).Value = 0.02795372904983374
$ws.Cells.Item(14, 17).Value = 0.21084051338
$ws.Cells.Item(14, 18).Value = 1.89756462042
$ws.Cells.Item(14, 19).Value = 0.0004185272024899925
$ws.Cells.Item(14, 20).Value = 0.0004185272024899926

# Row 15 updates
$ws.Cells.Item(15, 1).Value = "MuSCs"
$ws.Cells.Item(15, 5).Value = 2.0
$ws.Cells.Item(15, 6).Value = 0.6666666666666666
$ws.Cells.Item(15, 7).Value = 0.11543
$ws.Cells.Item(15, 8).Value = 0.34629
$ws.Cells.Item(15, 9).Value = 0.01497214206175765
$ws.Cells.Item(15, 10).Value = 0.01497214206175765
$ws.Cells.Item(15, 13).Value = 44.29005966666667
$ws.Cells.Item(15, 15).Value = 0.6778141756295529
$ws.Cells.Item(15, 16).Value = 0.6778141756295529
$ws.Cells.Item(15, 17).Value = 5.112401587323333
$ws.Cells.Item(15, 18).Value = 46.01161428591
$ws.Cells.Item(15, 19).Value = 0.01014833012899882
$ws.Cells.Item(15, 20).Value = 0.01014833012899882

# Row 16 updates
$ws.Cells.Item(16, 1).Value = "MuSCs"
$ws.Cells.Item(16, 5).Value = 2.0
$ws.Cells.Item(16, 6).Value = 0.6666666666666666
$ws.Cells.Item(16, 7).Value = 0.11543
$ws.Cells.Item(16, 8).Value = 0.34629
$ws.Cells.Item(16, 9).Value = 0.01497214206175765
$ws.Cells.Item(16, 10).Value = 0.01497214206175765
$ws.Cells.Item(16, 13).Value = 19.10886933333333
$ws.Cells.Item(16, 14).Value = 57.326608
$ws.Cells.Item(16, 15).Value = 0.2924417490485847
$ws.Cells.Item(16, 16).Value = 0.2924417490485847
$ws.Cells.Item(16, 17).Value = 2.205736787146666
$ws.Cells.Item(16, 18).Value = 19.85163108432
$ws.Cells.Item(16, 19).Value = 0.00437847941154429
$ws.Cells.Item(16, 20).Value = 0.004378479411544291

# Row 17 (new)
$ws.Cells.Item(17, 1).Value = "MuSCs"
$ws.Cells.Item(17, 2).Value = "Tnfsf13"
$ws.Cells.Item(17, 3).Value = "Sdc2"
$ws.Cells.Item(17, 4).Value = "Resolving-Mac"
$ws.Cells.Item(17, 5).Value = 2.0
$ws.Cells.Item(17, 6).Value = 0.6666666666666666
$ws.Cells.Item(17, 7).Value = 0.11543
$ws.Cells.Item(17, 8).Value = 0.34629
$ws.Cells.Item(17, 9).Value = 0.01497214206175765
$ws.Cells.Item(17, 10).Value = 0.01497214206175765
$ws.Cells.Item(17, 11).Value = 2.0
$ws.Cells.Item(17, 12).Value = 0.6666666666666666
$ws.Cells.Item(17, 13).Value = 0.1169856666666667
$ws.Cells.Item(17, 14).Value = 0.350957
$ws.Cells.Item(17, 15).Value = 0.001790346272028586
$ws.Cells.Item(17, 16).Value = 0.001790346272028586
$ws.Cells.Item(17, 17).Value = 0.01350365550333333
$ws.Cells.Item(17, 18).Value = 0.12153289953
$ws.Cells.Item(17, 19).Value = 0.0000268053187245502
$ws.Cells.Item(17, 20).Value = 0.00002680531872455021

# Row 18 (new)
$ws.Cells.Item(18, 1).Value = "Resolving-Mac"
$ws.Cells.Item(18, 2).Value = "Tnfsf13"
$ws.Cells.Item(18, 3).Value = "Sdc2"
$ws.Cells.Item(18, 4).Value = "ECs"
$ws.Cells.Item(18, 5).Value = 3.0
$ws.Cells.Item(18, 6).Value = 1.0
$ws.Cells.Item(18, 7).Value = 3.685254666666667
$ws.Cells.Item(18, 8).Value = 11.055764
$ws.Cells.Item(18, 9).Value = 0.4780053400596784
$ws.Cells.Item(18, 10).Value = 0.4780053400596784
$ws.Cells.Item(18, 11).Value = 3.0
$ws.Cells.Item(18, 12).Value = 1.0
$ws.Cells.Item(18, 13).Value = 1.826566
$ws.Cells.Item(18, 14).Value = 5.479698
$ws.Cells.Item(18, 15).Value = 0.02795372904983374
$ws.Cells.Item(18, 16).Value = 0.02795372904983374
$ws.Cells.Item(18, 17).Value = 6.731360875474667
$ws.Cells.Item(18, 18).Value = 60.582247879272
$ws.Cells.Item(18, 19).Value = 0.01336203176040189
$ws.Cells.Item(18, 20).Value = 0.01336203176040189

# Row 19 (new)
$ws.Cells.Item(19, 1).Value = "Resolving-Mac"
$ws.Cells.Item(19, 2).Value = "Tnfsf13"
$ws.Cells.Item(19, 3).Value = "Sdc2"
$ws.Cells.Item(19, 4).Value = "FAPs"
$ws.Cells.Item(19, 5).Value = 3.0
$ws.Cells.Item(19, 6).Value = 1.0
$ws.Cells.Item(19, 7).Value = 3.685254666666667
$ws.Cells.Item(19, 8).Value = 11.055764
$ws.Cells.Item(19, 9).Value = 0.4780053400596784
$ws.Cells.Item(19, 10).Value = 0.4780053400596784
$ws.Cells.Item(19, 11).Value = 3.0
$ws.Cells.Item(19, 12).Value = 1.0
$ws.Cells.Item(19, 13).Value = 44.29005966666667
$ws.Cells.Item(19, 14).Value = 132.870179
$ws.Cells.Item(19, 15).Value = 0.6778141756295529
$ws.Cells.Item(19, 16).Value = 0.6778141756295529
$ws.Cells.Item(19, 17).Value = 163.2201490735285
$ws.Cells.Item(19, 18).Value = 1468.981341661756
$ws.Cells.Item(19, 19).Value = 0.323998795519075
$ws.Cells.Item(19, 20).Value = 0.3239987955190751

# Row 20 (new)
$ws.Cells.Item(20, 1).Value = "Resolving-Mac"
$ws.Cells.Item(20, 2).Value = "Tnfsf13"
$ws.Cells.Item(20, 3).Value = "Sdc2"
$ws.Cells.Item(20, 4).Value = "MuSCs"
$ws.Cells.Item(20, 5).Value = 3.0
$ws.Cells.Item(20, 6).Value = 1.0
$ws.Cells.Item(20, 7).Value = 3.685254666666667
$ws.Cells.Item(20, 8).Value = 11.055764
$ws.Cells.Item(20, 9).Value = 0.4780053400596784
$ws.Cells.Item(20, 10).Value = 0.4780053400596784
$ws.Cells.Item(20, 11).Value = 3.0
$ws.Cells.Item(20, 12).Value = 1.0
$ws.Cells.Item(20, 13).Value = 19.10886933333333
$ws.Cells.Item(20, 14).Value = 57.326608
$ws.Cells.Item(20, 15).Value = 0.2924417490485847
$ws.Cells.Item(20, 16).Value = 0.2924417490485847
$ws.Cells.Item(20, 17).Value = 70.42104988539023
$ws.Cells.Item(20, 18).Value = 633.789448968512
$ws.Cells.Item(20, 19).Value = 0.1397887177016159
$ws.Cells.Item(20, 20).Value = 0.1397887177016159

# Row 21 (new)
$ws.Cells.Item(21, 1).Value = "Resolving-Mac"
$ws.Cells.Item(21, 2).Value = "Tnfsf13"
$ws.Cells.Item(21, 3).Value = "Sdc2"
$ws.Cells.Item(21, 4).Value = "Resolving-Mac"
$ws.Cells.Item(21, 5).Value = 3.0
$ws.Cells.Item(21, 6).Value = 1.0
$ws.Cells.Item(21, 7).Value = 3.685254666666667
$ws.Cells.Item(21, 8).Value = 11.055764
$ws.Cells.Item(21, 9).Value = 0.4780053400596784
$ws.Cells.Item(21, 10).Value = 0.4780053400596784
$ws.Cells.Item(21, 11).Value = 2.0
$ws.Cells.Item(21, 12).Value = 0.6666666666666666
$ws.Cells.Item(21, 13).Value = 0.1169856666666667
$ws.Cells.Item(21, 14).Value = 0.350957
$ws.Cells.Item(21, 15).Value = 0.001790346272028586
$ws.Cells.Item(21, 16).Value = 0.001790346272028586
$ws.Cells.Item(21, 17).Value = 0.4311219740164444
$ws.Cells.Item(21, 18).Value = 3.880097766147999
$ws.Cells.Item(21, 19).Value = 0.0008557950785856018
$ws.Cells.Item(21, 20).Value = 0.000855795078585602
